$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values
$ws.Range("A2").Value = 320.4131556305247
$ws.Range("B2").Value = 3289.743
$ws.Range("C2").Value = -775.769

# Add new row 3
$ws.Range("A3").Value = 441.6047391640358
$ws.Range("B3").Value = 3289.743
$ws.Range("C3").Value = -775.769

# Add new row 4
$ws.Range("A4").Value = 431.9462277151515
$ws.Range("B4").Value = 3289.743
$ws.Range("C4").Value = -775.769
